$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that Word leaves behind after the
#    last edit location. Deleting it causes the remaining "_Hlk..." bookmark
#    to be renumbered down automatically (id 1 -> id 0), matching the target.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. "Achievements:" -> "Achievement:" -- but only the SECOND occurrence
#    (the "Reduced backup time by 40% ... operation." bullet). Locate the
#    first occurrence, then restrict the search range to the remainder of
#    the document before replacing.
$first = $d.Content
$first.Find.Execute("Achievements:") | Out-Null
$rest = $d.Range($first.End, $d.Content.End)
$rest.Find.Execute("Achievements:", $true, $false, $false, $false, $false, $true, 1, $false, "Achievement:", 2) | Out-Null

# 3. "operations." -> "operation." in the same bullet (unique in the doc).
$d.Content.Find.Execute("operations.", $true, $false, $false, $false, $false, $true, 1, $false, "operation.", 2) | Out-Null

# 4. Expand the "Backend API Development" bullet description.
$oldBullet = " Built secure APIs with JWT & OAuth2, optimized DB queries, and integrated H2 for rapid testing."
$newBullet = " Developed secure APIs with JWT & OAuth2, managed DB with JPA/Hibernate, integrated H2 for rapid testing, handled global exceptions, and used RestTemplate and WebClient for service interaction."
$d.Content.Find.Execute($oldBullet, $true, $false, $false, $false, $false, $true, 1, $false, $newBullet, 2) | Out-Null
